$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start clean so the shared-strings table only contains the strings we
# actually need (old UI_COMPONENT_COUNT / int / How many UI components
# strings are dropped entirely in the target file).
$ws.Cells.Clear()

# Write cells in an order chosen so that newly-introduced unique strings
# appear in this sequence (matching the target workbook's shared string
# table order):
#   0 FILENAME, 1 POSITION_X, 2 POSITION_Y, 3 floor.png, 4 bullet.png,
#   5 SIZE_X, 6 SIZE_Y, 7 (size note)
$note = "이미지 원래 사이즈를 사용하고 싶으면 0"

$ws.Range("A1").Value = "FILENAME"
$ws.Range("A2").Value = "POSITION_X"
$ws.Range("A3").Value = "POSITION_Y"
$ws.Range("B1").Value = "floor.png"
$ws.Range("B6").Value = "bullet.png"
$ws.Range("A4").Value = "SIZE_X"
$ws.Range("A5").Value = "SIZE_Y"
$ws.Range("C4").Value = $note

$ws.Range("B2").Value = 0
$ws.Range("B3").Value = 0
$ws.Range("B4").Value = 0
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = $note

$ws.Range("A6").Value = "FILENAME"
$ws.Range("A7").Value = "POSITION_X"
$ws.Range("B7").Value = 100
$ws.Range("A8").Value = "POSITION_Y"
$ws.Range("B8").Value = 100
$ws.Range("A9").Value = "SIZE_X"
$ws.Range("B9").Value = 200
$ws.Range("C9").Value = $note
$ws.Range("A10").Value = "SIZE_Y"
$ws.Range("B10").Value = 200
$ws.Range("C10").Value = $note

# Column widths: stored widths are snapped to a 1/7-character pixel grid
# by this engine, so feed values that land on the desired grid point.
$ws.Columns.Item(1).ColumnWidth = 25.3   # -> stored width 26
$ws.Columns.Item(2).ColumnWidth = 11.7   # -> stored width ~12.43 (closest to 12.375)

$ws.Range("B10").Select() | Out-Null
